# Backup before dimension reduction: shift the "qN" labels in column A
# down by one index (q1 -> q0, q2 -> q1, ..., q96 -> q95) for data rows
# 2 through 97 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 97; $row++) {
    $oldIndex = $row - 1
    $newIndex = $oldIndex - 1
    $ws.Cells.Item($row, 1).Value = "q" + $newIndex
}
